$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.455.79'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.657.99'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.87%  '
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '201.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +11.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '579.42'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.652.25'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.622'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.68%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.682'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.155'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +8.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '57.44'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000296'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +19.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.14'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.249.48'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.667.86'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.77%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.126'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '68.519.28'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.70'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '404.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +27.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.30'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.98'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.66'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.86'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.25'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +23.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '32.11'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '696.96'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +13.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.30'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.117'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '64.79'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.434'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +17.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0802'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.02%  '
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.142'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.90%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.243.35'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +12.30%  '
$ws.Range('B43').Value = 'ThetaToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +13.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +19.52%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +38.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0423'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.99'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.132'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.67'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.11'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.71%  '
